$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the "Ciclo: ..." paragraph (currently the last paragraph in the doc).
$cicloPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ciclo:*") {
        $cicloPara = $p
    }
}

# 1. Replace the whole paragraph's contents so "Ciclo: 2016A" becomes two runs:
#    "Ciclo: 2016" and "ª" (same run formatting), mirroring how Word splits a
#    run when the trailing character is produced by a separate edit/autocorrect.
$rPr = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='32'/></w:rPr>"
$cicloXml = "<w:p $wNs w:rsidR='003E6B31' w:rsidRPr='003E6B31' w:rsidRDefault='003E6B31'><w:pPr>$rPr</w:pPr>" +
            "<w:r>$rPr<w:t>Ciclo: 2016</w:t></w:r>" +
            "<w:r>$rPr<w:t>&#170;</w:t></w:r>" +
            "</w:p>"
$cicloPara.Range.InsertXML($cicloXml)

# 2. Append a blank paragraph followed by a "Materia: ..." paragraph at the
#    very end of the document, matching the formatting used throughout.
$tailXml = "<w:p $wNs><w:pPr>$rPr</w:pPr></w:p>" +
           "<w:p $wNs><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>Materia: Programaci&#243;n WEB</w:t></w:r></w:p>"

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($tailXml)
